# Adds one new week (Fecha = 44522) of price data for "Ají" at the top of the
# existing data block: five new rows are inserted at row 494, pushing the
# previously-existing rows 494:561 down to 499:566.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows before the current row 494, shifting everything below down.
$ws.Rows("494:498").Insert()

# Columns that are constant across every data row in this sheet.
$ws.Range("A494:A498").Value = 6
$ws.Range("B494:B498").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C494:C498").Value = "Metropolitana"
$ws.Range("D494:D498").Value = 44522
$ws.Range("E494:E498").Value = 13
$ws.Range("F494:F498").Value = 100112021
$ws.Range("G494:G498").Value = "Ají"
$ws.Range("R494:R498").Value = "Hortaliza"

# Row 494: Ají, Americana (o), Primera
$ws.Range("H494").Value = "Americana (o)"
$ws.Range("I494").Value = "Primera"
$ws.Range("J494").Value = 80
$ws.Range("K494").Value = 20000
$ws.Range("L494").Value = 25000
$ws.Range("M494").Value = 23125
$ws.Range("N494").Value = "$/caja 25 kilos"
$ws.Range("O494").Value = "Provincia de Limarí"
$ws.Range("P494").Value = 925
$ws.Range("Q494").Value = 25

# Row 495: Ají, Americana (o), Segunda
$ws.Range("H495").Value = "Americana (o)"
$ws.Range("I495").Value = "Segunda"
$ws.Range("J495").Value = 15
$ws.Range("K495").Value = 15000
$ws.Range("L495").Value = 15000
$ws.Range("M495").Value = 15000
$ws.Range("N495").Value = "$/caja 25 kilos"
$ws.Range("O495").Value = "Provincia de Limarí"
$ws.Range("P495").Value = 600
$ws.Range("Q495").Value = 25

# Row 496: Ají, Chilena(o), Primera
$ws.Range("H496").Value = "Chilena(o)"
$ws.Range("I496").Value = "Primera"
$ws.Range("J496").Value = 35
$ws.Range("K496").Value = 55000
$ws.Range("L496").Value = 60000
$ws.Range("M496").Value = 57857
$ws.Range("N496").Value = "$/caja 25 kilos"
$ws.Range("O496").Value = "Provincia de Huasco"
$ws.Range("P496").Value = 2314
$ws.Range("Q496").Value = 25

# Row 497: Ají, Americana (o), Segunda
$ws.Range("H497").Value = "Americana (o)"
$ws.Range("I497").Value = "Segunda"
$ws.Range("J497").Value = 10
$ws.Range("K497").Value = 50000
$ws.Range("L497").Value = 50000
$ws.Range("M497").Value = 50000
$ws.Range("N497").Value = "$/caja 25 kilos"
$ws.Range("O497").Value = "Provincia de Huasco"
$ws.Range("P497").Value = 2000
$ws.Range("Q497").Value = 25

# Row 498: Ají, Inferno, Primera
$ws.Range("H498").Value = "Inferno"
$ws.Range("I498").Value = "Primera"
$ws.Range("J498").Value = 130
$ws.Range("K498").Value = 18000
$ws.Range("L498").Value = 20000
$ws.Range("M498").Value = 18769
$ws.Range("N498").Value = "$/caja 12 kilos"
$ws.Range("O498").Value = "Región de Arica y Parinacota"
$ws.Range("P498").Value = 1564
$ws.Range("Q498").Value = 12
